# ============================================================================
# Clinic management system.docx - "version two" edit
# Applies the changes described by the commit diff using Word COM automation.
# ============================================================================

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: apply the "List Paragraph" style to a paragraph while preserving
# the direct sz/szCs (12pt) run/paragraph-mark formatting that the style
# switch otherwise clears.
# ---------------------------------------------------------------------------
function Set-ListParagraphStyle($para) {
    $para.Range.Font.Size = 12
    $para.Format.Style = "List Paragraph"
    $para.Range.Font.Size = 12
    $para.Range.Font.SizeBi = 12
}

# ---------------------------------------------------------------------------
# STEP 1: "Use cases: patient walks in" paragraph (#17)
#   - becomes ListParagraph styled
#   - text changes
#   - gains a large block of new bullet paragraphs right after it
# ---------------------------------------------------------------------------

$pUse = $d.Paragraphs.Item(17)
Set-ListParagraphStyle $pUse
$pUse.Range.Text = "Use case: patient uses system to make an appointment"

# Insert the new block of paragraphs right after paragraph 17 (before its
# trailing paragraph mark) using one multi-line InsertAfter call.
$pUse = $d.Paragraphs.Item(17)
$insPos = $pUse.Range.End - 1
$insRange = $d.Range($insPos, $insPos)
$newBlock = "`rPatient logs in to the system and create an account or go straight to appointments" + `
    "`r-patient fills in general information:" + `
    "`r -full name" + `
    "`r-gender" + `
    "`r-age" + `
    "`r-phone number" + `
    "`r-appointment date" + `
    "`r-appointment time(optional)" + `
    "`r-specialty " + `
    "`r-doctor"
$insRange.InsertAfter($newBlock)

# Paragraph 18: "Patient logs in to the system ..." - ListParagraph + new
# numbered list (numId 7, freshly minted abstract numbering definition).
$p18 = $d.Paragraphs.Item(18)
Set-ListParagraphStyle $p18
$p18.Range.ListFormat.ApplyListTemplate($null)

# Paragraph 19: "-patient fills in general information:" - ListParagraph,
# left indent 1080 twips = 54pt.
$p19 = $d.Paragraphs.Item(19)
Set-ListParagraphStyle $p19
$p19.Format.LeftIndent = 54

# Paragraphs 20-27: the "-full name" ... "-doctor" sub bullets, all
# ListParagraph with left indent 1440 twips = 72pt.
for ($i = 20; $i -le 27; $i++) {
    $p = $d.Paragraphs.Item($i)
    Set-ListParagraphStyle $p
    $p.Format.LeftIndent = 72
}

Write-Output "Step1 paragraph count: $($d.Paragraphs.Count)"
for ($i = 16; $i -le 29; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Output "$i : [$($p.Range.Text)]"
}
